$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update project #8 (row 9): add description text, project name, and flip status to "done"
$ws.Range("D9").Value = "This project is a LIRI (Language Interpretation and Recognition Interface); a command line node app that takes in parameters and gives you back data."
$ws.Range("E9").Value = "liri-node-app"
$ws.Range("F9").Value = "done"

# Row 9 grows to accommodate the wrapped description text
$ws.Rows.Item(9).RowHeight = 45

# Move the active selection to the updated cell
$ws.Range("H9").Select()
